$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '28.876.58'
Set-TextValue 'D3' '1.825.27'
Set-TextValue 'D4' '0.9995'
Set-TextValue 'E4' '  -0.14%  '
Set-TextValue 'D5' '240.21'
Set-TextValue 'E5' '  -1.31%  '
Set-TextValue 'D6' '0.6905'
Set-TextValue 'E6' '  -2.33%  '
Set-TextValue 'D7' '0.9998'
Set-TextValue 'E7' '  -0.12%  '
Set-TextValue 'D8' '0.07610'
Set-TextValue 'E8' '  -3.16%  '
Set-TextValue 'D9' '0.3019'
Set-TextValue 'E9' '  -3.87%  '
Set-TextValue 'D10' '23.41'
Set-TextValue 'E10' '  -4.36%  '
Set-TextValue 'D11' '0.07744'
Set-TextValue 'E11' '  -3.48%  '
Set-TextValue 'D12' '1.829.07'
Set-TextValue 'E12' '  -2.10%  '
Set-TextValue 'D13' '5.040'
Set-TextValue 'E13' '  -3.12%  '
Set-TextValue 'D14' '90.07'
Set-TextValue 'E14' '  -3.57%  '
Set-TextValue 'D15' '0.6714'
Set-TextValue 'E15' '  -4.22%  '
Set-TextValue 'D16' '6.357'
Set-TextValue 'E16' '  -1.53%  '
Set-TextValue 'D17' '0.000008253'
Set-TextValue 'E17' '  -1.42%  '
Set-TextValue 'D18' '28.859.68'
Set-TextValue 'E18' '  -1.99%  '
Set-TextValue 'D19' '242.27'
Set-TextValue 'E19' '  -4.07%  '
Set-TextValue 'D20' '2.075.58'
Set-TextValue 'E20' '  -2.31%  '
Set-TextValue 'D21' '12.61'
Set-TextValue 'E21' '  -4.06%  '
Set-TextValue 'D22' '0.9996'
Set-TextValue 'E22' '  -0.11%  '
Set-TextValue 'D23' '7.395'
Set-TextValue 'E23' '  -2.77%  '
Set-TextValue 'D24' '0.9995'
Set-TextValue 'E24' '  -0.16%  '
Set-TextValue 'E25' '  -5.71%  '
Set-TextValue 'D26' '160.91'
Set-TextValue 'E26' '  +0.05%  '
Set-TextValue 'D27' '8.711'
Set-TextValue 'E27' '  -3.38%  '
Set-TextValue 'E28' '  -3.17%  '
Set-TextValue 'D29' '1.530'
Set-TextValue 'E29' '  +1.99%  '
Set-TextValue 'D30' '4.188'
Set-TextValue 'E30' '  -3.14%  '
Set-TextValue 'D31' '4.119'
Set-TextValue 'E31' '  -3.84%  '
Set-TextValue 'D32' '1.193'
Set-TextValue 'E32' '  -1.45%  '
Set-TextValue 'D33' '0.05082'
Set-TextValue 'E33' '  -4.24%  '
Set-TextValue 'D34' '0.7524'
Set-TextValue 'E34' '  +0.03%  '
Set-TextValue 'D35' '1.806'
Set-TextValue 'E35' '  -4.18%  '
Set-TextValue 'D36' '1.135'
Set-TextValue 'E36' '  -2.72%  '
Set-TextValue 'D37' '2.678'
Set-TextValue 'D38' '0.01836'
Set-TextValue 'E38' '  -2.30%  '
Set-TextValue 'D39' '1.196.23'
Set-TextValue 'E39' '  -5.46%  '
Set-TextValue 'D40' '2.676'
Set-TextValue 'E40' '  -2.36%  '
Set-TextValue 'D41' '0.9293'
Set-TextValue 'E41' '  +3.54%  '
Set-TextValue 'D42' '108.00'
Set-TextValue 'E42' '  -1.13%  '
Set-TextValue 'D43' '0.9991'
Set-TextValue 'D44' '0.5161'
Set-TextValue 'E44' '  -0.30%  '
Set-TextValue 'D45' '1.975.51'
Set-TextValue 'E45' '  -2.80%  '
Set-TextValue 'B46' 'EnergySwap'
Set-TextValue 'C46' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D46' '9.446'
Set-TextValue 'E46' '  -0.95%  '
Set-TextValue 'B47' 'BabyDogeCoin'
Set-TextValue 'C47' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D47' '0.00000000121'
Set-TextValue 'E47' '  -6.66%  '
Set-TextValue 'D48' '1.721'
Set-TextValue 'E48' '  -3.85%  '
Set-TextValue 'B49' 'FraxShare'
Set-TextValue 'C49' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D49' '5.181'
Set-TextValue 'E49' '  -13.16%  '
Set-TextValue 'B50' 'Aave'
Set-TextValue 'C50' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D50' '62.03'
Set-TextValue 'E50' '  -13.08%  '
Set-TextValue 'D51' '6.860'
Set-TextValue 'E51' '  -2.97%  '
